$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen the hourly data columns (D:P) -----------------------------------
# Target stored column width (OOXML) is 10.5703125 characters. The COM
# ColumnWidth setter here adds a fixed padding of 5/6 (0.8333...) character
# widths on write-back, so back that out of the desired stored width.
$ws.Columns("D:P").ColumnWidth = 10.5703125 - 0.8333333333333334

# --- Row 5 ("Beaming" / North America row) gets band values added ----------
# Numeric "20" for the outer columns, textual "20 15" for the middle band.
$ws.Range("D5:F5").Value = 20
$ws.Range("G5:J5").Value = "20 15"
$ws.Range("K5:N5").Value = 20

# --- Data grid (D5:P15) drops back to the default (non-bold, 11pt) font ----
$ws.Range("D5:P15").Font.Size = 11

# --- Cursor moves to S9 -------------------------------------------------
$null = $ws.Range("S9").Select()
